$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.168.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.764.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.764.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.43%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("E12").Value = "  -2.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000268"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.393.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.767.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.208.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.65%  "
$ws.Range("E19").Value = "  -3.85%  "
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "469.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  -5.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("E25").Value = "  -5.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.908.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("E31").Value = "  -4.43%  "
$ws.Range("E32").Value = "  -3.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.719.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("E38").Value = "  -3.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.77%  "
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "399.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.04%  "
